# The "CreatedUser" column (I) is blank for rows 3-7 while every other row
# already mirrors the "Username" column (F) into column I (e.g. rows 2, 8-16).
# Fill in the missing values the same way, copying straight from column F so
# the result always matches whatever is currently in the Username column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 6).Value2
}
